$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new "Price" values look like plain decimal numbers (e.g. "1.001" or
# "217.66"). Setting .Value directly on a General-formatted cell would let Excel
# auto-convert them to numbers, whereas the source data keeps them as text.
# Mark those specific cells as Text before writing the value so they round-trip
# exactly like the original (multi-dot / subscript) price strings already do.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.653.97'
$ws.Range("E2").Value = '  +2.87%  '
$ws.Range("D3").Value = '1.688.91'
$ws.Range("E3").Value = '  +3.30%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '217.66'
$ws.Range("E5").Value = '  +3.90%  '
$ws.Range("D6").Value = '0.5335'
$ws.Range("E6").Value = '  +2.53%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").Value = '0.2682'
$ws.Range("E8").Value = '  +4.58%  '
$ws.Range("D9").Value = '0.06435'
$ws.Range("E9").Value = '  +3.07%  '
$ws.Range("D10").Value = '21.68'
$ws.Range("E10").Value = '  +6.83%  '
$ws.Range("D11").Value = '0.07805'
$ws.Range("E11").Value = '  +3.38%  '
$ws.Range("D12").Value = '1.695.79'
$ws.Range("E12").Value = '  +3.93%  '
$ws.Range("D13").Value = '4.504'
$ws.Range("E13").Value = '  +3.41%  '
$ws.Range("D14").Value = '0.5636'
$ws.Range("E14").Value = '  +4.08%  '
$ws.Range("D15").Value = '0.0₅8456'
$ws.Range("E15").Value = '  +6.71%  '
$ws.Range("D16").Value = '66.42'
$ws.Range("E16").Value = '  +3.04%  '
$ws.Range("D17").Value = '26.702.81'
$ws.Range("E17").Value = '  +3.05%  '
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").Value = '4.817'
$ws.Range("E19").Value = '  +4.18%  '
$ws.Range("D20").Value = '196.04'
$ws.Range("E20").Value = '  +6.43%  '
$ws.Range("D21").Value = '10.42'
$ws.Range("E21").Value = '  +4.21%  '
$ws.Range("D22").Value = '6.389'
$ws.Range("E22").Value = '  +5.19%  '
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("D24").Value = '144.16'
$ws.Range("E24").Value = '  -0.94%  '
$ws.Range("D25").Value = '0.1290'
$ws.Range("E25").Value = '  +7.61%  '
$ws.Range("D26").Value = '7.494'
$ws.Range("E26").Value = '  +2.18%  '
$ws.Range("D27").Value = '16.29'
$ws.Range("E27").Value = '  +5.26%  '
$ws.Range("D28").Value = '1.415'
$ws.Range("E28").Value = '  +2.97%  '
$ws.Range("D29").Value = '0.06188'
$ws.Range("E29").Value = '  +4.08%  '
$ws.Range("D30").Value = '1.281'
$ws.Range("E30").Value = '  +3.21%  '
$ws.Range("D31").Value = '3.605'
$ws.Range("E31").Value = '  +7.89%  '
$ws.Range("D32").Value = '3.473'
$ws.Range("E32").Value = '  +3.75%  '
$ws.Range("D33").Value = '1.707'
$ws.Range("E33").Value = '  +6.42%  '
$ws.Range("D34").Value = '1.016'
$ws.Range("E34").Value = '  +4.73%  '
$ws.Range("D35").Value = '2.801'
$ws.Range("E35").Value = '  +2.56%  '
$ws.Range("D36").Value = '2.422'
$ws.Range("E36").Value = '  +1.68%  '
$ws.Range("D37").Value = '0.5747'
$ws.Range("E37").Value = '  -0.98%  '
$ws.Range("E38").Value = '  +4.06%  '
$ws.Range("D39").Value = '6.017'
$ws.Range("E39").Value = '  +6.38%  '
$ws.Range("D40").Value = '1.078.22'
$ws.Range("E40").Value = '  +5.20%  '
$ws.Range("D41").Value = '0.8670'
$ws.Range("E41").Value = '  +3.31%  '
$ws.Range("D42").Value = '1.000'
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").Value = '100.53'
$ws.Range("E43").Value = '  +0.89%  '
$ws.Range("D44").Value = '1.839.81'
$ws.Range("E44").Value = '  +3.10%  '
$ws.Range("D45").Value = '0.0₈109'
$ws.Range("E45").Value = '  +3.07%  '
$ws.Range("D46").Value = '57.46'
$ws.Range("E46").Value = '  +5.91%  '
$ws.Range("D47").Value = '8.211'
$ws.Range("E47").Value = '  +3.07%  '
$ws.Range("D48").Value = '1.002'
$ws.Range("E48").Value = '  +0.13%  '
$ws.Range("D49").Value = '0.05224'
$ws.Range("E49").Value = '  +0.88%  '
$ws.Range("D50").Value = '6.101'
$ws.Range("E50").Value = '  +5.40%  '
$ws.Range("D51").Value = '0.4241'
$ws.Range("E51").Value = '  +0.37%  '
